# Renames the three inline logo pictures living in the document's
# headers/footers:
#   - Footer (first page),  docPr id="3"  Pearson logo: image1.png -> image2.png
#   - Footer (default),     docPr id="2"  Pearson logo: image1.png -> image2.png
#   - Header (first page),  docPr id="1"  BTec logo:    image2.jpg -> image1.jpg
#
# InlineShapes don't expose a settable Name through Range.InlineShapes
# directly (Word reports the handle as stale), so each picture is selected
# first and then addressed through $word.Selection.InlineShapes, which is
# the path Word actually resolves writes through.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Footers: the Pearson Edexcel logo ------------------------------------
# wdHeaderFooterPrimary = 1 (the "default" footer/footer2.xml, docPr id="2")
# wdHeaderFooterFirstPage = 2 (the "first page" footer/footer1.xml, docPr id="3")
$firstPageFooterShape = $sec.Footers.Item(2).Range.InlineShapes.Item(1)
$firstPageFooterShape.Range.Select()
$word.Selection.InlineShapes.Item(1).Name = "image2.png"

$primaryFooterShape = $sec.Footers.Item(1).Range.InlineShapes.Item(1)
$primaryFooterShape.Range.Select()
$word.Selection.InlineShapes.Item(1).Name = "image2.png"

# --- Header: the BTec logo --------------------------------------------------
# wdHeaderFooterFirstPage = 2 (the "first page" header/header1.xml, docPr id="1")
$firstPageHeaderShape = $sec.Headers.Item(2).Range.InlineShapes.Item(1)
$firstPageHeaderShape.Range.Select()
$word.Selection.InlineShapes.Item(1).Name = "image1.jpg"

Write-Output "Renamed inline picture names in footers and header."
